$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F (想去人数) values for several rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 136
$wsExhibit.Range("F5").Value = 372
$wsExhibit.Range("F8").Value = 2085
$wsExhibit.Range("F9").Value = 10790
$wsExhibit.Range("F12").Value = 291
$wsExhibit.Range("F14").Value = 423
$wsExhibit.Range("F15").Value = 9012
$wsExhibit.Range("F16").Value = 1118
$wsExhibit.Range("F18").Value = 5289
$wsExhibit.Range("F20").Value = 3356

# Sheet "全部类型" (sheet4) - update column F (想去人数) values for several rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 136
$wsAll.Range("F5").Value = 372
$wsAll.Range("F9").Value = 2085
$wsAll.Range("F12").Value = 10790
$wsAll.Range("F15").Value = 291
$wsAll.Range("F17").Value = 423
$wsAll.Range("F18").Value = 9012
$wsAll.Range("F19").Value = 1118
$wsAll.Range("F21").Value = 5289
$wsAll.Range("F23").Value = 3356
